$wb = $excel.ActiveWorkbook

# Sheet "展览" (全部 exhibit data): F3 1271 -> 1274, F4 2766 -> 2772
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1274
$wsExhibit.Range("F4").Value = 2772

# Sheet "全部类型" (all types aggregate): F5 1271 -> 1274, F6 2766 -> 2772
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1274
$wsAll.Range("F6").Value = 2772
